$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 20:55:00"
$wsZhCn.Range("H2").Value = "2016-03-18 20:55:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 20:55:04"
$wsDeDe.Range("H2").Value = "2016-03-18 20:55:36"
